# Refresh cached market-board figures (currentAveragePrice / NQ / HQ) and the
# derived Leve price + profit columns (H:N) for the affected leve rows across
# all eight crafting-job sheets, per the scheduled-runner market data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70: Consecrating Congregation
$ws.Range("H70").Value = 54306.9
$ws.Range("I70").Value = 1215.3334
$ws.Range("K70").Value = 3646.0002
$ws.Range("M70").Value = -3376.0002

# Row 73: Curbing the Contagion (L)
$ws.Range("H73").Value = 54306.9
$ws.Range("I73").Value = 1215.3334
$ws.Range("K73").Value = 3646.0002
$ws.Range("M73").Value = -2710.0002

# Row 76: Warding Off Temptation
$ws.Range("H76").Value = 43481132
$ws.Range("I76").Value = 47621904
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 47621904
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -47621589
$ws.Range("N76").Value = -3630

# Row 79: The Garden of Arcane Delights (L)
$ws.Range("H79").Value = 43481132
$ws.Range("I79").Value = 47621904
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 47621904
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -47620812
$ws.Range("N79").Value = -5184

# Row 80: Cleansing the Wicked Humours
$ws.Range("H80").Value = 12512708
$ws.Range("I80").Value = 30550
$ws.Range("J80").Value = 20834148
$ws.Range("K80").Value = 91650
$ws.Range("L80").Value = 62502444
$ws.Range("M80").Value = -90652
$ws.Range("N80").Value = -62504440

# Row 83: Washing Away the Sins (L)
$ws.Range("H83").Value = 12512708
$ws.Range("I83").Value = 30550
$ws.Range("J83").Value = 20834148
$ws.Range("K83").Value = 274950
$ws.Range("L83").Value = 187507332
$ws.Range("M83").Value = -269958
$ws.Range("N83").Value = -187517316

# Row 86: Filling in the Blanks
$ws.Range("H86").Value = 633097.6
$ws.Range("I86").Value = 900
$ws.Range("J86").Value = 843830.2
$ws.Range("K86").Value = 900
$ws.Range("L86").Value = 843830.2
$ws.Range("M86").Value = 223
$ws.Range("N86").Value = -846076.2

# Row 89: Ink into Antiquity (L)
$ws.Range("H89").Value = 633097.6
$ws.Range("I89").Value = 900
$ws.Range("J89").Value = 843830.2
$ws.Range("K89").Value = 4500
$ws.Range("L89").Value = 4219151
$ws.Range("M89").Value = 1116
$ws.Range("N89").Value = -4230383

# Row 92: Whinier than the Sword
$ws.Range("H92").Value = 564.3929000000001
$ws.Range("I92").Value = 477.23077
$ws.Range("J92").Value = 1697.5
$ws.Range("K92").Value = 477.23077
$ws.Range("L92").Value = 1697.5
$ws.Range("M92").Value = 770.76923
$ws.Range("N92").Value = -4193.5

# Row 98: The Dotted Line
$ws.Range("H98").Value = 39529580
$ws.Range("I98").Value = 16668493
$ws.Range("K98").Value = 16668493
$ws.Range("M98").Value = -16666995

# Row 103: Let Loose the Juice
$ws.Range("H103").Value = 50002260
$ws.Range("I103").Value = 111112260
$ws.Range("J103").Value = 3172.7273
$ws.Range("K103").Value = 333336780
$ws.Range("L103").Value = 9518.1819
$ws.Range("M103").Value = -333336194
$ws.Range("N103").Value = -10690.1819

# Row 106: Making Your Mark
$ws.Range("H106").Value = 35716188
$ws.Range("I106").Value = 50001404
$ws.Range("J106").Value = 3150
$ws.Range("K106").Value = 50001404
$ws.Range("L106").Value = 3150
$ws.Range("M106").Value = -50000773
$ws.Range("N106").Value = -4412

# Row 116: Growing Up
$ws.Range("H116").Value = 33346190
$ws.Range("I116").Value = 22224222
$ws.Range("J116").Value = 41687668
$ws.Range("K116").Value = 22224222
$ws.Range("L116").Value = 41687668
$ws.Range("M116").Value = -22220780
$ws.Range("N116").Value = -41694552

# Row 122: Wishful Inking
$ws.Range("H122").Value = 39529580
$ws.Range("I122").Value = 16668493
$ws.Range("K122").Value = 50005479
$ws.Range("M122").Value = -50003029

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 28866156
$ws.Range("I137").Value = 6580086
$ws.Range("K137").Value = 19740258
$ws.Range("M137").Value = -19737708

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 8125190.5
$ws.Range("I32").Value = 1960932.8
$ws.Range("J32").Value = 38472308
$ws.Range("K32").Value = 1960932.8
$ws.Range("L32").Value = 38472308
$ws.Range("M32").Value = -1960645.8
$ws.Range("N32").Value = -38472882

# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 417952.84
$ws.Range("I45").Value = 501262.7
$ws.Range("J45").Value = 1403.5
$ws.Range("K45").Value = 501262.7
$ws.Range("L45").Value = 1403.5
$ws.Range("M45").Value = -500885.7
$ws.Range("N45").Value = -2157.5

# Row 97: Ore for Me
$ws.Range("H97").Value = 874.1818
$ws.Range("I97").Value = 824
$ws.Range("J97").Value = 1100
$ws.Range("K97").Value = 824
$ws.Range("L97").Value = 1100
$ws.Range("M97").Value = -328
$ws.Range("N97").Value = -2092

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 24313272
$ws.Range("I132").Value = 25649754
$ws.Range("J132").Value = 18521852
$ws.Range("K132").Value = 76949262
$ws.Range("L132").Value = 55565556
$ws.Range("M132").Value = -76946732
$ws.Range("N132").Value = -55570616

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 1950.52
$ws.Range("I86").Value = 1958.1123
$ws.Range("J86").Value = 1578.5
$ws.Range("K86").Value = 1958.1123
$ws.Range("L86").Value = 1578.5
$ws.Range("M86").Value = -835.1123
$ws.Range("N86").Value = -3824.5

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 1950.52
$ws.Range("I89").Value = 1958.1123
$ws.Range("J89").Value = 1578.5
$ws.Range("K89").Value = 9790.5615
$ws.Range("L89").Value = 7892.5
$ws.Range("M89").Value = -4174.5615
$ws.Range("N89").Value = -19124.5

# Row 94: High Steal
$ws.Range("H94").Value = 1383.2667
$ws.Range("I94").Value = 949.9231
$ws.Range("J94").Value = 4200
$ws.Range("K94").Value = 949.9231
$ws.Range("L94").Value = 4200
$ws.Range("M94").Value = -498.9231
$ws.Range("N94").Value = -5102

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 2472.6667
$ws.Range("I105").Value = 2419.25
$ws.Range("K105").Value = 2419.25
$ws.Range("M105").Value = -672.25

$ws = $wb.Worksheets.Item("CRP")
# Row 86: Birch, Please
$ws.Range("H86").Value = 6601.241
$ws.Range("I86").Value = 8921.5
$ws.Range("J86").Value = 3745.5386
$ws.Range("K86").Value = 8921.5
$ws.Range("L86").Value = 3745.5386
$ws.Range("M86").Value = -7798.5
$ws.Range("N86").Value = -5991.5386

# Row 89: Built This City on Blocks and Soul (L)
$ws.Range("H89").Value = 6601.241
$ws.Range("I89").Value = 8921.5
$ws.Range("J89").Value = 3745.5386
$ws.Range("K89").Value = 44607.5
$ws.Range("L89").Value = 18727.693
$ws.Range("M89").Value = -38991.5
$ws.Range("N89").Value = -29959.693

# Row 105: Zelkova, My Love
$ws.Range("H105").Value = 8091.1875
$ws.Range("I105").Value = 2006
$ws.Range("J105").Value = 14176.375
$ws.Range("K105").Value = 2006
$ws.Range("L105").Value = 14176.375
$ws.Range("M105").Value = -259
$ws.Range("N105").Value = -17670.375

# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 8675.272000000001
$ws.Range("I122").Value = 10314.223
$ws.Range("K122").Value = 30942.669
$ws.Range("M122").Value = -28492.669

$ws = $wb.Worksheets.Item("CUL")
# Row 132: More Mezcal
$ws.Range("H132").Value = 3411.6667
$ws.Range("I132").Value = 3251.1667
$ws.Range("K132").Value = 29260.5003
$ws.Range("M132").Value = -26730.5003

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit
$ws.Range("H70").Value = 2330620.5
$ws.Range("I70").Value = 1467462.1
$ws.Range("J70").Value = 3642621.2
$ws.Range("K70").Value = 1467462.1
$ws.Range("L70").Value = 3642621.2
$ws.Range("M70").Value = -1467192.1
$ws.Range("N70").Value = -3643161.2

# Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 2330620.5
$ws.Range("I73").Value = 1467462.1
$ws.Range("J73").Value = 3642621.2
$ws.Range("K73").Value = 1467462.1
$ws.Range("L73").Value = 3642621.2
$ws.Range("M73").Value = -1466526.1
$ws.Range("N73").Value = -3644493.2

# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 10035.5
$ws.Range("I80").Value = 4969
$ws.Range("J80").Value = 12968.737
$ws.Range("K80").Value = 4969
$ws.Range("L80").Value = 12968.737
$ws.Range("M80").Value = -3971
$ws.Range("N80").Value = -14964.737

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 10035.5
$ws.Range("I83").Value = 4969
$ws.Range("J83").Value = 12968.737
$ws.Range("K83").Value = 24845
$ws.Range("L83").Value = 64843.685
$ws.Range("M83").Value = -19853
$ws.Range("N83").Value = -74827.685

# Row 126: Gold Rush Order
$ws.Range("H126").Value = 6124.973
$ws.Range("I126").Value = 21442
$ws.Range("J126").Value = 3731.6875
$ws.Range("K126").Value = 64326
$ws.Range("L126").Value = 11195.0625
$ws.Range("M126").Value = -61856
$ws.Range("N126").Value = -16135.0625

$ws = $wb.Worksheets.Item("LTW")
# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 1854744.2
$ws.Range("I132").Value = 2565380.2
$ws.Range("K132").Value = 7696140.600000001
$ws.Range("M132").Value = -7693610.600000001

$ws = $wb.Worksheets.Item("WVR")
# Row 135: In Line with Linen
$ws.Range("H135").Value = 45886
$ws.Range("J135").Value = 45886
$ws.Range("L135").Value = 45886
$ws.Range("N135").Value = -56026

Write-Output "Updated market pricing for $($wb.Worksheets.Count) sheets (215 cells)"
